$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 11:52"

# Row 9 - Alemania
$ws.Range("B9").Value = 159934
$ws.Range("C9").Value = 22
$ws.Range("E9").Value = 33220

# Row 51 - Malasia
$ws.Range("B51").Value = 5945
$ws.Range("C51").Value = 94
$ws.Range("D51").Value = 4087
$ws.Range("E51").Value = 1758
$ws.Range("F51").Value = 40

# Row 54 - Finlandia
$ws.Range("B54").Value = 4906
$ws.Range("C54").Value = 166
$ws.Range("E54").Value = 1907

# Row 68 - Uzbekistan
$ws.Range("D68").Value = 1005
$ws.Range("E68").Value = 942

# Row 83 - Eslovenia
$ws.Range("B83").Value = 1418
$ws.Range("C83").Value = 10
$ws.Range("D83").Value = 230
$ws.Range("F83").Value = 25
$ws.Range("G83").Value = 3
$ws.Range("H83").Value = 89
